# edit.ps1
# Applies the commit's changes to Taxonomias.docx:
#  1. Delete the list item "Actualizaciones regulares de firmware: ..."
#     (its text no longer appears anywhere in the "good security features" list).
#  2. Re-run-split (with w:proofErr spell-check markers) the paragraphs whose
#     text is unchanged but whose run/proofErr structure changed as a side
#     effect of Word's live spell-checker re-validating the edited area:
#       - "Autenticación fuerte de usuario: ..."
#       - "Cifrado de datos: ..."            (now adjacent to item 1 above)
#       - "Contraseñas predeterminadas débiles: ..."
#  3. Move the two <w:lastRenderedPageBreak/> markers down to the next
#     list item, reflecting the reflow caused by the deleted paragraph.

$d = $word.ActiveDocument

$wordMlNs = 'xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"'

function Get-FoundRange($doc, [string]$text) {
    $probe = $doc.Content
    $ok = $probe.Find.Execute($text, $false, $false, $false, $false, $false, $true, 1, $false, "", 0)
    if (-not $ok) {
        throw "Text not found: $text"
    }
    return $doc.Range($probe.Start, $probe.End)
}

function Set-RangeBodyXml($doc, [string]$text, [string]$bodyInnerXml) {
    $target = Get-FoundRange $doc $text
    $pkg = '<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage">' +
           '<pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml">' +
           '<pkg:xmlData><w:document ' + $wordMlNs + '><w:body>' + $bodyInnerXml + '</w:body></w:document></pkg:xmlData>' +
           '</pkg:part></pkg:package>'
    $target.InsertXML($pkg)
}

# ---------------------------------------------------------------------------
# 1. Delete the whole "Actualizaciones regulares de firmware..." list item.
# ---------------------------------------------------------------------------
$oldFirmwareText = "Actualizaciones regulares de firmware: Actualizaciones regulares para corregir vulnerabilidades de seguridad y problemas conocidos."
$rng = $d.Content
$ok = $rng.Find.Execute($oldFirmwareText, $false, $false, $false, $false, $false, $true, 1, $false, "", 0)
if (-not $ok) { throw "Could not find firmware paragraph" }

$paraToDelete = $null
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $p = $d.Paragraphs.Item($i)
    if ($p.Range.Start -le $rng.Start -and $p.Range.End -ge $rng.End) {
        $paraToDelete = $p
        break
    }
}
if ($null -eq $paraToDelete) { throw "Could not locate paragraph object to delete" }
$paraToDelete.Range.Delete()

# ---------------------------------------------------------------------------
# 2. Re-split "Autenticación fuerte de usuario..." with spell-check proofErr
#    markup (text content is unchanged).
# ---------------------------------------------------------------------------
$autenticacionText = "Autenticación fuerte de usuario: Requerir una autenticación sólida, como la autenticación de dos factores, para acceder al dispositivo."
$autenticacionBody = '<w:p>' +
    '<w:proofErr w:type="spellStart"/>' +
    '<w:r><w:t>Autenticación</w:t></w:r>' +
    '<w:proofErr w:type="spellEnd"/>' +
    '<w:r><w:t xml:space="preserve"> </w:t></w:r>' +
    '<w:proofErr w:type="spellStart"/>' +
    '<w:r><w:t>fuerte</w:t></w:r>' +
    '<w:proofErr w:type="spellEnd"/>' +
    '<w:r><w:t xml:space="preserve"> de </w:t></w:r>' +
    '<w:proofErr w:type="spellStart"/>' +
    '<w:r><w:t>usuario</w:t></w:r>' +
    '<w:proofErr w:type="spellEnd"/>' +
    '<w:r><w:t xml:space="preserve">: Requerir una autenticación sólida, como la autenticación de dos factores, para acceder al </w:t></w:r>' +
    '<w:proofErr w:type="spellStart"/>' +
    '<w:r><w:t>dispositivo</w:t></w:r>' +
    '<w:proofErr w:type="spellEnd"/>' +
    '<w:r><w:t>.</w:t></w:r>' +
    '</w:p>'
Set-RangeBodyXml $d $autenticacionText $autenticacionBody

# ---------------------------------------------------------------------------
# 3. Re-split "Cifrado de datos..." with a proofErr marker around "Cifrado"
#    (text content is unchanged, it just sits next to item 1 now).
# ---------------------------------------------------------------------------
$cifradoText = "Cifrado de datos: La encriptación de datos transmitidos y almacenados, lo que dificulta su acceso no autorizado."
$cifradoBody = '<w:p>' +
    '<w:proofErr w:type="spellStart"/>' +
    '<w:r><w:t>Cifrado</w:t></w:r>' +
    '<w:proofErr w:type="spellEnd"/>' +
    '<w:r><w:t xml:space="preserve"> de datos: La encriptación de datos transmitidos y almacenados, lo que dificulta su acceso no autorizado.</w:t></w:r>' +
    '</w:p>'
Set-RangeBodyXml $d $cifradoText $cifradoBody

# ---------------------------------------------------------------------------
# 4. Re-split "Contraseñas predeterminadas débiles..." with proofErr markers
#    (text content is unchanged).
# ---------------------------------------------------------------------------
$contrasenasText = "Contraseñas predeterminadas débiles: Contraseñas que son fáciles de adivinar o que no se pueden cambiar fácilmente."
$contrasenasBody = '<w:p>' +
    '<w:proofErr w:type="spellStart"/>' +
    '<w:r><w:t>Contraseñas</w:t></w:r>' +
    '<w:proofErr w:type="spellEnd"/>' +
    '<w:r><w:t xml:space="preserve"> </w:t></w:r>' +
    '<w:proofErr w:type="spellStart"/>' +
    '<w:r><w:t>predeterminadas</w:t></w:r>' +
    '<w:proofErr w:type="spellEnd"/>' +
    '<w:r><w:t xml:space="preserve"> </w:t></w:r>' +
    '<w:proofErr w:type="spellStart"/>' +
    '<w:r><w:t>débiles</w:t></w:r>' +
    '<w:proofErr w:type="spellEnd"/>' +
    '<w:r><w:t xml:space="preserve">: </w:t></w:r>' +
    '<w:proofErr w:type="spellStart"/>' +
    '<w:r><w:t>Contraseñas</w:t></w:r>' +
    '<w:proofErr w:type="spellEnd"/>' +
    '<w:r><w:t xml:space="preserve"> que son fáciles de adivinar o que no se pueden cambiar fácilmente.</w:t></w:r>' +
    '</w:p>'
Set-RangeBodyXml $d $contrasenasText $contrasenasBody

# ---------------------------------------------------------------------------
# 5. Move the first <w:lastRenderedPageBreak/> one list item further down:
#    from "Falta de autenticación de usuario..." to "Aplicaciones de
#    terceros inseguras...".
# ---------------------------------------------------------------------------
$faltaAutenticacionText = "Falta de autenticación de usuario: Permitir que cualquier persona pueda acceder y controlar el dispositivo sin autenticación adecuada."
$faltaAutenticacionBody = '<w:p><w:r><w:t>' + $faltaAutenticacionText + '</w:t></w:r></w:p>'
Set-RangeBodyXml $d $faltaAutenticacionText $faltaAutenticacionBody

$aplicacionesText = "Aplicaciones de terceros inseguras: La falta de seguridad en las aplicaciones de terceros que se utilizan en conjunción con el dispositivo."
$aplicacionesBody = '<w:p><w:r><w:lastRenderedPageBreak/><w:t>' + $aplicacionesText + '</w:t></w:r></w:p>'
Set-RangeBodyXml $d $aplicacionesText $aplicacionesBody

# ---------------------------------------------------------------------------
# 6. Move the second <w:lastRenderedPageBreak/> one list item further down:
#    from "Desperdicio de energía..." to "Vida útil corta...".
# ---------------------------------------------------------------------------
$desperdicioText = "Desperdicio de energía: Los dispositivos IoT que no tienen un modo de suspensión o que no pueden apagarse pueden desperdiciar energía innecesariamente."
$desperdicioBody = '<w:p><w:r><w:t>' + $desperdicioText + '</w:t></w:r></w:p>'
Set-RangeBodyXml $d $desperdicioText $desperdicioBody

$vidaUtilText = "Vida útil corta: Los dispositivos IoT que están diseñados para durar poco tiempo pueden aumentar la cantidad de residuos electrónicos."
$vidaUtilBody = '<w:p><w:r><w:lastRenderedPageBreak/><w:t>' + $vidaUtilText + '</w:t></w:r></w:p>'
Set-RangeBodyXml $d $vidaUtilText $vidaUtilBody

Write-Host "All edits applied."
